$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 401.5
$ws.Range("I4").Value = 401.5
$ws.Range("K4").Value = 401.5
$ws.Range("M4").Value = -287.5

$ws.Range("H33").Value = 179.375
$ws.Range("I33").Value = 162.5
$ws.Range("K33").Value = 162.5
$ws.Range("M33").Value = 66.5

$ws.Range("H40").Value = 1103.1428
$ws.Range("I40").Value = 1103.1428
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1103.1428
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -928.1428000000001

$ws.Range("H42").Value = 386.6
$ws.Range("I42").Value = 388.75
$ws.Range("J42").Value = 378
$ws.Range("K42").Value = 1166.25
$ws.Range("L42").Value = 1134
$ws.Range("M42").Value = -936.25
$ws.Range("N42").Value = -1594

$ws.Range("H64").Value = 4329.6665
$ws.Range("I64").Value = 4494.5
$ws.Range("K64").Value = 4494.5
$ws.Range("M64").Value = -4246.5

$ws.Range("H67").Value = 4329.6665
$ws.Range("I67").Value = 4494.5
$ws.Range("K67").Value = 4494.5
$ws.Range("M67").Value = -3636.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 118
$ws.Range("I5").Value = 118
$ws.Range("K5").Value = 118
$ws.Range("M5").Value = -6

$ws.Range("H74").Value = 2995
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 2995
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H88").Value = 1217.8
$ws.Range("I88").Value = 1575
$ws.Range("K88").Value = 1575
$ws.Range("M88").Value = -1169

$ws.Range("H91").Value = 1217.8
$ws.Range("I91").Value = 1575
$ws.Range("K91").Value = 1575
$ws.Range("M91").Value = -171

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 118
$ws.Range("I4").Value = 118
$ws.Range("K4").Value = 118
$ws.Range("M4").Value = -3

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0

$ws.Range("H86").Value = 2962.8
$ws.Range("I86").Value = 1523
$ws.Range("K86").Value = 1523
$ws.Range("M86").Value = -400

$ws.Range("H89").Value = 2962.8
$ws.Range("I89").Value = 1523
$ws.Range("K89").Value = 7615
$ws.Range("M89").Value = -1999

$ws.Range("H105").Value = 1754.1666
$ws.Range("I105").Value = 1550
$ws.Range("K105").Value = 1550
$ws.Range("M105").Value = 197

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 36.8
$ws.Range("I7").Value = 23.5
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 23.5
$ws.Range("L7").Value = 90
$ws.Range("M7").Value = 89.5
$ws.Range("N7").Value = -316

$ws.Range("H22").Value = 1449
$ws.Range("I22").Value = 749.3333
$ws.Range("K22").Value = 749.3333
$ws.Range("M22").Value = -399.3333

$ws.Range("H28").Value = 32601.334
$ws.Range("J28").Value = 34656.184
$ws.Range("L28").Value = 34656.184
$ws.Range("N28").Value = -35146.184

$ws.Range("H36").Value = 9487
$ws.Range("I36").Value = 9487
$ws.Range("K36").Value = 9487
$ws.Range("M36").Value = -9099

$ws.Range("H40").Value = 9487
$ws.Range("I40").Value = 9487
$ws.Range("K40").Value = 9487
$ws.Range("M40").Value = -9327

$ws.Range("H134").Value = 9999
$ws.Range("I134").Value = 9999
$ws.Range("K134").Value = 29997
$ws.Range("M134").Value = -27462

$ws.Range("H141").Value = 249735.12
$ws.Range("J141").Value = 249735.12
$ws.Range("L141").Value = 249735.12
$ws.Range("N141").Value = -260095.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.6667
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 2700
$ws.Range("N5").Value = -2924

$ws.Range("H35").Value = 375
$ws.Range("I35").Value = 375
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1125
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -837

$ws.Range("H42").Value = 1487
$ws.Range("J42").Value = 1487
$ws.Range("L42").Value = 4461
$ws.Range("N42").Value = -5529

$ws.Range("H43").Value = 2987.5
$ws.Range("J43").Value = 2987.5
$ws.Range("L43").Value = 8962.5
$ws.Range("N43").Value = -9190.5

$ws.Range("H54").Value = 2213.75
$ws.Range("J54").Value = 2213.75
$ws.Range("L54").Value = 6641.25
$ws.Range("N54").Value = -7759.25

$ws.Range("H59").Value = 15
$ws.Range("I59").Value = 15
$ws.Range("K59").Value = 45
$ws.Range("M59").Value = 495

$ws.Range("H98").Value = 500.55554
$ws.Range("I98").Value = 497
$ws.Range("J98").Value = 502.33334
$ws.Range("K98").Value = 1491
$ws.Range("L98").Value = 1507.00002
$ws.Range("M98").Value = 7
$ws.Range("N98").Value = -4503.000019999999

$ws.Range("H104").Value = 2088
$ws.Range("I104").Value = 2088
$ws.Range("K104").Value = 6264
$ws.Range("M104").Value = -3643

$ws.Range("H107").Value = 1195.5
$ws.Range("I107").Value = 930
$ws.Range("K107").Value = 2790
$ws.Range("M107").Value = -870

$ws.Range("H116").Value = 2699.75
$ws.Range("I116").Value = 933.3333
$ws.Range("K116").Value = 2799.9999
$ws.Range("M116").Value = 642.0001000000002

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws.Range("H124").Value = 1000
$ws.Range("I124").Value = 1000
$ws.Range("K124").Value = 3000
$ws.Range("M124").Value = 1910

$ws.Range("H135").Value = 799.6667
$ws.Range("J135").Value = 900
$ws.Range("L135").Value = 8100
$ws.Range("N135").Value = -13170

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws.Range("H140").Value = 415.2857
$ws.Range("I140").Value = 415.2857
$ws.Range("K140").Value = 1245.8571
$ws.Range("M140").Value = 3934.1429

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7999.3335

$ws.Range("H73").Value = 7999.3335

$ws.Range("H80").Value = 24243.125
$ws.Range("I80").Value = 18511.25
$ws.Range("K80").Value = 18511.25
$ws.Range("M80").Value = -17513.25

$ws.Range("H83").Value = 24243.125
$ws.Range("I83").Value = 18511.25
$ws.Range("K83").Value = 92556.25
$ws.Range("M83").Value = -87564.25

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0

$ws.Range("H136").Value = 24081.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 24081.5
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").Value = 72244.5
$ws.Range("N136").Value = -77344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1388

$ws.Range("H40").Value = 1499
$ws.Range("I40").Value = 1499
$ws.Range("K40").Value = 1499
$ws.Range("M40").Value = -1363

$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -2030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 33737.2
$ws.Range("I58").Value = 21499.666
$ws.Range("K58").Value = 21499.666
$ws.Range("M58").Value = -21191.666

$ws.Range("H81").Value = 900
$ws.Range("I81").Value = 900
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -739

$ws.Range("H84").Value = 900
$ws.Range("I84").Value = 900
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -3696

$ws.Range("H96").Value = 1985.6
$ws.Range("I96").Value = 1146.3334
$ws.Range("K96").Value = 1146.3334
$ws.Range("M96").Value = 226.6666

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0

$ws.Range("H136").Value = 2421.6667
$ws.Range("I136").Value = 1911.875
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 5735.625
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -3185.625
$ws.Range("N136").Value = -24600
